$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 172 (pushes existing rows 172..272 down to 175..275)
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()
$ws.Rows.Item(172).Insert()

# New row 172 - Sandia, Perú, Primera
$ws.Cells.Item(172,1).Value = 4
$ws.Cells.Item(172,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(172,3).Value = "Los Lagos"
$ws.Cells.Item(172,4).Value = 44875
$ws.Cells.Item(172,5).Value = 10
$ws.Cells.Item(172,6).Value = 100112028
$ws.Cells.Item(172,7).Value = "Sandia"
$ws.Cells.Item(172,8).Value = "Sin especificar"
$ws.Cells.Item(172,9).Value = "Primera"
$ws.Cells.Item(172,10).Value = 200
$ws.Cells.Item(172,11).Value = 1400
$ws.Cells.Item(172,12).Value = 1400
$ws.Cells.Item(172,13).Value = 1400
$ws.Cells.Item(172,14).Value = "`$/unidad"
$ws.Cells.Item(172,15).Value = "Perú"
$ws.Cells.Item(172,16).Value = 1400
$ws.Cells.Item(172,17).Value = 1
$ws.Cells.Item(172,18).Value = "Hortaliza"

# New row 173 - Sandia, Perú, Segunda
$ws.Cells.Item(173,1).Value = 4
$ws.Cells.Item(173,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(173,3).Value = "Los Lagos"
$ws.Cells.Item(173,4).Value = 44875
$ws.Cells.Item(173,5).Value = 10
$ws.Cells.Item(173,6).Value = 100112028
$ws.Cells.Item(173,7).Value = "Sandia"
$ws.Cells.Item(173,8).Value = "Sin especificar"
$ws.Cells.Item(173,9).Value = "Segunda"
$ws.Cells.Item(173,10).Value = 200
$ws.Cells.Item(173,11).Value = 1300
$ws.Cells.Item(173,12).Value = 1300
$ws.Cells.Item(173,13).Value = 1300
$ws.Cells.Item(173,14).Value = "`$/unidad"
$ws.Cells.Item(173,15).Value = "Perú"
$ws.Cells.Item(173,16).Value = 1300
$ws.Cells.Item(173,17).Value = 1
$ws.Cells.Item(173,18).Value = "Hortaliza"

# New row 174 - Sandia, Perú, Tercera
$ws.Cells.Item(174,1).Value = 4
$ws.Cells.Item(174,2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(174,3).Value = "Los Lagos"
$ws.Cells.Item(174,4).Value = 44875
$ws.Cells.Item(174,5).Value = 10
$ws.Cells.Item(174,6).Value = 100112028
$ws.Cells.Item(174,7).Value = "Sandia"
$ws.Cells.Item(174,8).Value = "Sin especificar"
$ws.Cells.Item(174,9).Value = "Tercera"
$ws.Cells.Item(174,10).Value = 200
$ws.Cells.Item(174,11).Value = 950
$ws.Cells.Item(174,12).Value = 950
$ws.Cells.Item(174,13).Value = 950
$ws.Cells.Item(174,14).Value = "`$/unidad"
$ws.Cells.Item(174,15).Value = "Perú"
$ws.Cells.Item(174,16).Value = 950
$ws.Cells.Item(174,17).Value = 1
$ws.Cells.Item(174,18).Value = "Hortaliza"
